$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - Analyse concurrentielle (installation et test d'applications)
$ws.Range("A18").Value = "Analyse"
$ws.Range("B18").Value = "Analyse concurentielle, j'installe plusieurs application de tracking fitness sur mon téléphone et je note les fonctionnalités que je trouve intéressantes ainsi que celles que je trouve moins inéressantes"
$ws.Range("C18").Value = 1.5
$ws.Range("D18").Value = "2019-02-13"

# Row 19 - Continuation de l'analyse concurrentielle
$ws.Range("A19").Value = "Analyse"
$ws.Range("B19").Value = "Continuation de l'analyse concurentielle, à partir des tests que j'ai faits sur les différents applications, je remplis un tableau comparatif des fonctionnalités intégrées par les applications"
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Value = "2019-02-13"

# Row 20 - Définition de l'audience de l'application
$ws.Range("A20").Value = "Analyse"
$ws.Range("B20").Value = "Définition de l'audience de l'application. Il s'agit ici de définir les personnes qui pouraient être amenés à utilisé l'application qui sera développer"
$ws.Range("C20").Value = "0..5"
$ws.Range("D20").Value = "2019-02-13"

# Preserve wrap-text styling consistent with the rest of the sheet
$ws.Range("A18:A20").WrapText = $true
$ws.Range("B18:B20").WrapText = $true

[void]$ws.Range("E20").Select()
